# Clean version of catalan-juvenile-recidivism data
# Rename a handful of column-name values in column B of Sheet1 to their
# "cleaned" equivalents.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B16").Value  = "V12_n_criminal_record"
$ws.Range("B144").Value = "V132_RECID2013_recid"
$ws.Range("B126").Value = "V115_RECID2015_recid"
$ws.Range("B34").Value  = "V30_program_start"
$ws.Range("B35").Value  = "V31_program_end"
$ws.Range("B33").Value  = "V29_program_duration"
$ws.Range("B31").Value  = "V27_program_duration_cat"
$ws.Range("B25").Value  = "V22_main_crime_comission_date"
$ws.Range("B32").Value  = "V28_days_from_crime_to_program"

# Reflect the final cursor/scroll position left behind by the edits.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("B32").Select()
